# Daily attendance processing - 2025-11-29 05:25:53
# Normalizes the "Recorded By" column (G) on the session analysis sheet.
# For each data row, if the recorded-by list has more than one entry and
# the first entry is not already "System", rotate the list left by one
# position (move the first entry to the end of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value

    if ($value -eq $null) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value -notmatch ",") { continue }

    $parts = $value -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -lt 2) { continue }
    if ($parts[0] -eq "System") { continue }

    $rotated = $parts[1..($parts.Length - 1)] + @($parts[0])
    $newValue = [string]::Join(", ", $rotated)

    $cell.Value = $newValue
}
